# "Added previous raw output" (reverted): the RawOutput query table previously
# included two extra rows (c101c5.txt, rows 4-5) that are removed here, which
# also shrinks the table/ExternalData range and changes the dependent
# MIN/VLOOKUP results on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("RawOutput")

# Remove the last two raw-output rows (the c101c5.txt runs) entirely, which
# shrinks the sheet/table extent from A1:F5 down to A1:F3.
$ws2.Range("A4:F5").EntireRow.Delete()

# The RawOutput sheet keeps a new selection on A4 (the first empty row).
$ws2.Activate()
$ws2.Range("A4").Select()

# Shrink the ExternalData_1 defined name so it matches the smaller table.
$nm = $wb.Names.Item("RawOutput!ExternalData_1")
$nm.RefersTo = "='RawOutput'!`$A`$1:`$F`$3"

# Sheet1 becomes the active/selected tab again.
$ws1.Activate()
